$wb = $excel.ActiveWorkbook
$wsCasos = $wb.Worksheets.Item("Casos de Uso")
$wsInstructivo = $wb.Worksheets.Item("Instructivo")

# --- Update text for the two previously-placeholder use cases (rows 9 and 10) ---
$wsCasos.Cells.Item(9, 3).Value = "El líder de evento deberá poder crear un comité para poder llevar el control de los organizadores de manera fácil e intuitiva.`n"
$wsCasos.Cells.Item(10, 3).Value = "El líder de evento y comité deberá poder modificar los miembros del comité con el fin de agregar o elimianar miembros de manera sencilla e intuitiva.`n"

# --- Update estimated time column (F) for several rows ---
$wsCasos.Cells.Item(9, 6).Value = 15
$wsCasos.Cells.Item(10, 6).Value = 15
$wsCasos.Cells.Item(11, 6).Value = 15
$wsCasos.Cells.Item(12, 6).Value = 15
$wsCasos.Cells.Item(17, 6).Value = 14
$wsCasos.Cells.Item(19, 6).Value = 14

# --- Update sheet view / selection (set Instructivo first so "Casos de Uso" ends up the active/selected tab) ---
$wsInstructivo.Range("C8").Select()
$wsCasos.Range("C11").Select()
